$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.826.68"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.73%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.631.73"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.96%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.997"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.41%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.76"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.02%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.522"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.12%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.997"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.38%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.30"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.18%  "

# Row 9
$ws.Range("E9").Value = "  -2.89%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0614"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.25%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0882"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.98%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.868.08"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.77%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.637.92"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.70%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.05"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.79%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.565"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.32%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.27"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.69%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "27.880.31"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.52%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "230.79"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.71%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0₃0723"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.15%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.54"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.91%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.997"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.36%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.36"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.70%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.36"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.38%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.08"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.98%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "154.60"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.51%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.96"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.41%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.64"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.92%  "

# Row 28
$ws.Range("E28").Value = "  -1.18%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.998"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.30%  "

# Row 30
$ws.Range("E30").Value = "  -0.82%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0482"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.70%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.40"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.63%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.406.80"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.70%  "

# Row 34
$ws.Range("E34").Value = "  -0.15%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.56"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.04%  "

# Row 36
$ws.Range("E36").Value = "  +9.22%  "

# Row 37
$ws.Range("E37").Value = "  +1.07%  "

# Row 38
$ws.Range("E38").Value = "  +0.30%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.562"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.61%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.869"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.54%  "

# Row 41
$ws.Range("E41").Value = "  +0.22%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.998"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.29%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "66.73"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.77%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.52"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.75%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.83"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.19%  "

# Row 46
$ws.Range("E46").Value = "  -1.19%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.775.24"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.97%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "87.79"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.56%  "

# Row 49
$ws.Range("B49").Value = "Algorand"
$ws.Range("C49").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.1000"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.41%  "

# Row 50
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0506"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.47%  "

# Row 51
$ws.Range("B51").Value = "BabyDogeCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0₇0982"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -6.33%  "

